$d = $word.ActiveDocument

# Correct the capitalization of the title ("Genérico" -> "genérico")
$d.Content.Find.Execute("Ejercicio Genérico M2A: Rellenar huecos", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ejercicio genérico M2A: Rellenar huecos", 2)

# Move the "_GoBack" bookmark (Word's "last edit location" marker) so that it
# now sits inside the title run, right after "Ejercicio g". Adding a bookmark
# named "_GoBack" automatically removes any pre-existing "_GoBack" bookmark
# elsewhere in the document, which is exactly the effect seen in the diff
# (the bookmark that used to sit further down the document disappears).
$titlePara = $d.Paragraphs(1).Range
$bmStart = $titlePara.Start + 11
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
